$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.000378370285034
$ws.Range("B1").Value = 3.112407445907593
$ws.Range("C1").Value = 6.713497161865234
$ws.Range("D1").Value = 1.907366991043091
$ws.Range("E1").Value = 1.338062644004822
